$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column D slightly (stored OOXML width 12 -> 13)
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666

# Update VENTA (D) and POR CUMPLIR (E) for the "OTROS" row
$ws.Range("D2").Value = 2942.59
$ws.Range("E2").Value = -2942.59

# Update TOTAL row (VENTA, POR CUMPLIR, CUMPLIMIENTO)
$ws.Range("D4").Value = 3191.07
$ws.Range("E4").Value = 14308.93
$ws.Range("F4").Value = 0.1823468571428571
